$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.897.21'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '2.961.23'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.13'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.18'
$ws.Range('E6').Value = '  +5.62%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('E8').Value = '  +0.83%  '
$ws.Range('D9').Value = '2.958.85'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.60'
$ws.Range('E10').Value = '  -5.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('E12').Value = '  +1.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.17'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '66.086.37'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '3.464.87'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.87'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '2.974.53'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.97'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.69'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.674'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.18'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.85'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.22'
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.16'
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.96'
$ws.Range('E28').Value = '  -9.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.07'
$ws.Range('E29').Value = '  +5.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.40'
$ws.Range('E30').Value = '  +10.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.58'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').Value = '0.0₃0990'
$ws.Range('E32').Value = '  -10.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.06'
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.70'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.37'
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.00'
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '43.67'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.299'
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.83'
$ws.Range('E42').Value = '  -8.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.119'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.33'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '386.44'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0353'
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('D47').Value = '2.712.07'
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.55'
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.106'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.13'
$ws.Range('E51').Value = '  +4.84%  '
